$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Truth Table")

# Apply an AutoFilter on column O ("OpASel", 15th column in A:T) showing only rows with value "01"
$ws.Range("A1:T38").AutoFilter(15, @("01"), 7)

# Update the active selection on the sheet to P1 to match the post-edit state
$ws.Range("P1").Select()
